$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
  $co = $ws.ChartObjects(1)
  Write-Output $co.Name
} catch {
  Write-Output ("error1: " + $_.Exception.Message)
}
